# 2021 Excess Mortality Update
# Updates the Figure 3b table with revised 2020 figures and adds 2021 Q2-Q4 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# --- Revisions to existing column-C values (2020 quarters + 2021 Q1) ---
$ws.Range("C2").Value = -0.037
$ws.Range("C12").Value = 0.099
$ws.Range("C20").Value = 0.372
$ws.Range("C21").Value = 0.381
$ws.Range("C22").Value = 0.275
$ws.Range("C23").Value = 0.65
$ws.Range("C25").Value = 0.147
$ws.Range("C26").Value = 0.373
$ws.Range("C27").Value = 0.484
$ws.Range("C28").Value = 0.313
$ws.Range("C29").Value = 0.867
$ws.Range("C30").Value = 0.482
$ws.Range("C31").Value = 0.089

# --- New rows: 2021 Q2, Q3, Q4 by Race/Ethnicity ---
$newRows = @(
    @("2021 Q2", "American Indian or Alaska Native", 0.183),
    @("2021 Q2", "Asian", 0.025),
    @("2021 Q2", "Black", 0.016),
    @("2021 Q2", "Latino", 0.048),
    @("2021 Q2", "Native Hawaiian and other Pacific Islander", 0.099),
    @("2021 Q2", "White", -0.042),
    @("2021 Q3", "American Indian or Alaska Native", 0.42),
    @("2021 Q3", "Asian", 0.138),
    @("2021 Q3", "Black", 0.189),
    @("2021 Q3", "Latino", 0.276),
    @("2021 Q3", "Native Hawaiian and other Pacific Islander", 0.412),
    @("2021 Q3", "White", 0.097),
    @("2021 Q4", "American Indian or Alaska Native", 0.445),
    @("2021 Q4", "Asian", 0.099),
    @("2021 Q4", "Black", 0.105),
    @("2021 Q4", "Latino", 0.235),
    @("2021 Q4", "Native Hawaiian and other Pacific Islander", 0.215),
    @("2021 Q4", "White", 0.063)
)

$startRow = 32
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# --- Resize the table to cover the new rows ---
$lastRow = $startRow + $newRows.Count - 1
$lo.Resize($ws.Range("A1:C" + $lastRow))
